$p = $ppt.ActivePresentation

function Get-SlideByTitle($pres, $titleText) {
    for ($i = 1; $i -le $pres.Slides.Count; $i++) {
        $slide = $pres.Slides.Item($i)
        try {
            $t = $slide.Shapes.Item("Title 1").TextFrame.TextRange.Text
        } catch {
            $t = $null
        }
        if ($t -eq $titleText) {
            return $slide
        }
    }
    return $null
}

# --- 1. Insert a new "Feature selection" slide right before "Modeling approach" ---
# Duplicating "Modeling approach" gives the new slide the same "Title and Content"
# layout/formatting as the rest of the deck. Duplicate() drops the copy immediately
# after its source, so move it one slot earlier to land right before the original.
$modeling = Get-SlideByTitle $p "Modeling approach"
$modelingIndex = $modeling.SlideIndex

$dupRange = $modeling.Duplicate()
$newSlide = $dupRange.Item(1)
$newSlide.MoveTo($modelingIndex)

$newSlide.Shapes.Item("Title 1").TextFrame.TextRange.Text = "Feature selection"
$newSlide.Shapes.Item("Content Placeholder 2").TextFrame.TextRange.Text = "Avoided variables with signs of survivorship bias"

# --- 2. Fill in the (previously empty) body placeholder of the "Results" slide ---
$results = Get-SlideByTitle $p "Results"
$results.Shapes.Item("Content Placeholder 2").TextFrame.TextRange.Text = "Predicting damage: no substantial improvement over baseline"
